# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the content is
# now "Ready for handoff" (previously "In Translation"), and refreshes the
# handoff timestamps. Also widens the Status column(s) to fit the new,
# longer status text (mirrors Excel's column autosize behaviour).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps refreshed as part of the new handoff generation -----------
$overview.Range("G2").Value = "2016-08-18 16:38:09"
$zhcn.Range("H2").Value     = "2016-08-18 16:37:58"
$dede.Range("H2").Value     = "2016-08-18 16:38:09"

# --- Widen the Status column(s) so the longer text fits (autofit) ---------
# (target raw OOXML column width ~17.216 chars; the host quantizes
# ColumnWidth in 1/6-char steps, so 16.3333 is the input that lands on the
# nearest achievable width.)
$overview.Columns.Item(5).ColumnWidth = 16.3333
$overview.Columns.Item(6).ColumnWidth = 16.3333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333
$dede.Columns.Item(3).ColumnWidth     = 16.3333
